# Updated Masterdata as per 2nd may Data Refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where regcntr_id (column A) needs to change to 10003
$rowsToUpdate = @(3, 23, 43, 63, 83, 105, 114, 123, 132, 141)
foreach ($r in $rowsToUpdate) {
    $ws.Cells.Item($r, 1).Value2 = 10003
}

# Update the view's selection to the row right after the data (select
# full rows from 162 to the end, as the last saved selection)
$ws.Range("A162:XFD1048576").Select()
